$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds one new daily price record for "Bruselas (repollito)" at
# Vega Modelo de Temuco. In the source data rows are kept in (roughly)
# reverse-chronological insertion order, so the new record is inserted as
# row 47, pushing the former rows 47-143 down to 48-144.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new record's data. Columns
# A, B, C, E, F, G, H, I, N, O, Q, R repeat the same values used throughout
# this sheet/product block.
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 44838
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100112035
$ws.Cells.Item(47, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 15
$ws.Cells.Item(47, 11).Value = 24000
$ws.Cells.Item(47, 12).Value = 24000
$ws.Cells.Item(47, 13).Value = 24000
$ws.Cells.Item(47, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(47, 16).Value = 2400
$ws.Cells.Item(47, 17).Value = 10
$ws.Cells.Item(47, 18).Value = "Hortaliza"
